# Bold the deck title on slide 1 ("Constrains in wheat production")
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Title
$title.TextFrame.TextRange.Font.Bold = $true
